# Update countries & provincias Spain
# Applies the "1 de Mayo de 2020 a las 11:52" data refresh to the Pais sheet:
#  - updates the "last updated" timestamp cell
#  - several countries swap relative ranking (rows shift) bringing new
#    country names + new daily figures into rows that previously held a
#    different country.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 11:52"

# Helper-less direct row rewrites: Country, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(35,  "Emiratos Arabes Unidos",       13038, 557, 2543, 10384, 1, 6, 111),
    @(36,  "Polonia",                      12877,   0, 3236,  8997, 160, 0, 644),
    @(54,  "Finlandia",                     5051,  56, 3000,  1840,  48, 0, 211),
    @(137, "Birmania",                       151,   1,   28,   117,   0, 0,   6),
    @(165, "Republica de Africa Central",     64,  14,   10,    54,   0, 0,   0),
    @(166, "Benin",                           64,   0,   33,    30,   0, 0,   1),
    @(167, "Libia",                           61,   0,   18,    40,   0, 0,   3),
    @(168, "Nepal",                           59,   2,   16,    43,   0, 0,   0),
    @(169, "Polinesia Francesa",              58,   0,   50,     8,   1, 0,   0),
    @(191, "Santo Tome y Principe",           16,   2,    4,    11,   0, 1,   1),
    @(192, "San Vicente y las Granadinas",    16,   0,    8,     8,   0, 0,   0),
    @(193, "Namibia",                         16,   0,    8,     8,   0, 0,   0),
    @(194, "Dominica",                        16,   0,   13,     3,   0, 0,   0),
    @(195, "Curazao",                         16,   0,   13,     2,   0, 0,   1),
    @(196, "Tayikistan",                      15,   0,    0,    15,   0, 0,   0),
    @(197, "San Cristobal y Nieves",          15,   0,    6,     9,   0, 0,   0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
